# Update 19-Mei-2021, end of day update.
# Petty cash book: fill in Buku KAS HARIAN ("Sheet1") daily transactions for
# 17/18/19-Mei-2021 (rows 3-22), then move the view/selection forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 (17-Mei-2021, SALDO AWAL carry row) ---
$ws.Range("D3").Value2 = 0

# --- Row 4: TRANSFER BCA ---
$ws.Range("B4").Value2 = "TRANSFER BCA"
$ws.Range("D4").Formula = "=1004000+6470000+3800000+25600000+10080000"

# --- Row 5: PLN - Astar 165 ---
$ws.Range("B5").Value2 = "PLN - Astar 165"
$ws.Range("D5").Formula = "=895000"

# --- Row 6: TELPON 5224823 ---
$ws.Range("B6").Value2 = "TELPON 5224823"
$ws.Range("D6").Formula = "=223500"

# --- Row 7: A/R ---
$ws.Range("B7").Value2 = "A/R"
$ws.Range("C7").Formula = "=6470000+25600000+10080000+10834000"

# --- Row 8: SALES - cash/retail ---
$ws.Range("B8").Value2 = "SALES - cash/retail"
$ws.Range("C8").Formula = "=5275475+21002525-10834000"

# --- Row 9: PRIVE - andreas ---
$ws.Range("B9").Value2 = "PRIVE - andreas"
$ws.Range("D9").Value2 = 2000000

# --- Row 10: SELISIH - lebih ---
$ws.Range("B10").Value2 = "SELISIH - lebih"
$ws.Range("C10").Formula = "=50000"

# --- Row 11: SETOR KE BANK ---
$ws.Range("B11").Value2 = "SETOR KE BANK"
$ws.Range("D11").Value2 = 18000000

# --- Row 12 (18-Mei-2021, Wages Expense) ---
$ws.Range("A12").Value2 = 44334
$ws.Range("B12").Value2 = "Wages Expense"
$ws.Range("D12").Formula = "=45000+195000"

# --- Row 13: A/R ---
$ws.Range("B13").Value2 = "A/R"
$ws.Range("C13").Formula = "=27020000+42656000"

# --- Row 14: TRANSFER BCA ---
$ws.Range("B14").Value2 = "TRANSFER BCA"
$ws.Range("D14").Formula = "=7000000+13633000+7508000+510000+497000+300000+1138000"

# --- Row 15: SALES - cash/retail ---
$ws.Range("B15").Value2 = "SALES - cash/retail"
$ws.Range("C15").Formula = "=2558475+49852525-42656000"

# --- Row 16: SELISIH - lebih ---
$ws.Range("B16").Value2 = "SELISIH - lebih"
$ws.Range("C16").Value2 = 18000

# --- Row 17: SETOR KE BANK ---
$ws.Range("B17").Value2 = "SETOR KE BANK"
$ws.Range("D17").Formula = "=49000000"

# --- Row 18 (19-Mei-2021, Wages Expense) ---
$ws.Range("A18").Value2 = 44335
$ws.Range("B18").Value2 = "Wages Expense"
$ws.Range("D18").Formula = "=45000"

# --- Row 19: TRANSFER BCA ---
$ws.Range("B19").Value2 = "TRANSFER BCA"
$ws.Range("D19").Formula = "=3640000+1250000+1170000+440000+515000"

# --- Row 20: A/R ---
$ws.Range("B20").Value2 = "A/R"
$ws.Range("C20").Formula = "=9300000"

# --- Row 21: TRANSFER BCA AA ---
$ws.Range("B21").Value2 = "TRANSFER BCA AA"
$ws.Range("D21").Formula = "=9300000"

# --- Row 22: LPG (note trailing space in the label) ---
$ws.Range("B22").Value2 = "LPG "
$ws.Range("D22").Value2 = 290000

# --- Move the frozen-pane scroll position / active selection to reflect
#     where the user ended up after today's entries (end-of-day update) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A17").Select()
$win.FreezePanes = $true
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("D38").Select()
